$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Expand the outlined/hidden Week 11 block (rows 49-55) and the
# detail columns (D,F,H,J,L,N,P) *before* writing any values into them -
# writing into a still-hidden row causes the engine to stamp a stray
# auto-fit row height, so unhide first.
$ws.Range("A49:A55").EntireRow.Hidden = $false
$ws.Range("D1:Q1").EntireColumn.Hidden = $false

# --- New shared strings must be created in this exact order so that the
# shared-string table indices line up with the target workbook:
#   95 = "Roos en Danial kwartier te laat"
#   96 = "goed doorgewerkt"
#   97 = "kwartier te laat"
# We therefore touch Q50 first, then the "goed doorgewerkt" cells, then the
# remaining "kwartier te laat" cells.

$ws.Range("Q50").Value = "Roos en Danial kwartier te laat"

# Week 10 (rows 41-47) - Woensdag (44) and Donderdag (45) remark columns
$ws.Range("D44").Value = "goed doorgewerkt"
$ws.Range("F44").Value = "goed doorgewerkt"
$ws.Range("H44").Value = "goed doorgewerkt"
$ws.Range("J44").Value = "goed doorgewerkt"
$ws.Range("L44").Value = "goed doorgewerkt"

$ws.Range("D45").Value = "goed doorgewerkt"
$ws.Range("F45").Value = "goed doorgewerkt"
$ws.Range("H45").Value = "goed doorgewerkt"
$ws.Range("J45").Value = "goed doorgewerkt"
$ws.Range("L45").Value = "goed doorgewerkt"
$ws.Range("N45").Value = "goed doorgewerkt"

# Week 10 - Vrijdag (46) totals row: hours + remarks
$ws.Range("B46").Value = 4
$ws.Range("C46").Value = 4
$ws.Range("D46").Value = "goed doorgewerkt"
$ws.Range("E46").Value = 4
$ws.Range("F46").Value = "goed doorgewerkt"
$ws.Range("F46").Borders.Item(9).LineStyle = -4142
$ws.Range("G46").Value = 4
$ws.Range("H46").Value = "goed doorgewerkt"
$ws.Range("H46").Borders.Item(9).LineStyle = -4142
$ws.Range("I46").Value = 4
$ws.Range("J46").Value = "goed doorgewerkt"
$ws.Range("J46").Borders.Item(9).LineStyle = -4142
$ws.Range("K46").Value = 4
$ws.Range("L46").Value = "goed doorgewerkt"
$ws.Range("L46").Borders.Item(9).LineStyle = -4142
$ws.Range("M46").Value = 4
$ws.Range("N46").Value = "goed doorgewerkt"
$ws.Range("N46").Borders.Item(9).LineStyle = -4142
$ws.Range("P46").Value = "Ziek"
$ws.Range("Q46").Value = "Danial Ziek"

# Week 11 (rows 49-55) - Maandag (50)
$ws.Range("B50").Value = 4
$ws.Range("C50").Value = 4
$ws.Range("D50").Value = "kwartier te laat"
$ws.Range("E50").Value = 4
$ws.Range("G50").Value = 4
$ws.Range("I50").Value = 4
$ws.Range("K50").Value = 4
$ws.Range("M50").Value = 4
$ws.Range("O50").Value = 4
$ws.Range("P50").Value = "kwartier te laat"

# Restore the view/selection state recorded after the edit.
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Range("D51").Select()
